$d = $word.ActiveDocument

$d.Content.Find.Execute("92×95=8740", $true, $false, $false, $false, $false, $true, 1, $false, "76×22=1672", 2)
$d.Content.Find.Execute("45×28=1260", $true, $false, $false, $false, $false, $true, 1, $false, "99×37=3663", 2)
$d.Content.Find.Execute("11×60=660", $true, $false, $false, $false, $false, $true, 1, $false, "27×30=810", 2)
$d.Content.Find.Execute("98×79=7742", $true, $false, $false, $false, $false, $true, 1, $false, "95×78=7410", 2)
$d.Content.Find.Execute("45×38=1710", $true, $false, $false, $false, $false, $true, 1, $false, "62×86=5332", 2)
$d.Content.Find.Execute("57×58=3306", $true, $false, $false, $false, $false, $true, 1, $false, "83×32=2656", 2)
$d.Content.Find.Execute("40×38=1520", $true, $false, $false, $false, $false, $true, 1, $false, "11×72=792", 2)
$d.Content.Find.Execute("37×89=3293", $true, $false, $false, $false, $false, $true, 1, $false, "79×26=2054", 2)
$d.Content.Find.Execute("54×49=2646", $true, $false, $false, $false, $false, $true, 1, $false, "68×86=5848", 2)
$d.Content.Find.Execute("79×96=7584", $true, $false, $false, $false, $false, $true, 1, $false, "16×28=448", 2)
$d.Content.Find.Execute("80×66=5280", $true, $false, $false, $false, $false, $true, 1, $false, "59×80=4720", 2)
$d.Content.Find.Execute("66×83=5478", $true, $false, $false, $false, $false, $true, 1, $false, "41×89=3649", 2)
$d.Content.Find.Execute("11×46=506", $true, $false, $false, $false, $false, $true, 1, $false, "73×62=4526", 2)
$d.Content.Find.Execute("59×59=3481", $true, $false, $false, $false, $false, $true, 1, $false, "34×64=2176", 2)
$d.Content.Find.Execute("61×24=1464", $true, $false, $false, $false, $false, $true, 1, $false, "38×56=2128", 2)
$d.Content.Find.Execute("53×67=3551", $true, $false, $false, $false, $false, $true, 1, $false, "54×63=3402", 2)
$d.Content.Find.Execute("35×48=1680", $true, $false, $false, $false, $false, $true, 1, $false, "35×61=2135", 2)
$d.Content.Find.Execute("25×93=2325", $true, $false, $false, $false, $false, $true, 1, $false, "33×30=990", 2)
$d.Content.Find.Execute("44×94=4136", $true, $false, $false, $false, $false, $true, 1, $false, "77×68=5236", 2)
$d.Content.Find.Execute("92×37=3404", $true, $false, $false, $false, $false, $true, 1, $false, "95×36=3420", 2)
$d.Content.Find.Execute("14×78=1092", $true, $false, $false, $false, $false, $true, 1, $false, "25×66=1650", 2)
$d.Content.Find.Execute("38×83=3154", $true, $false, $false, $false, $false, $true, 1, $false, "72×45=3240", 2)
$d.Content.Find.Execute("63×51=3213", $true, $false, $false, $false, $false, $true, 1, $false, "42×87=3654", 2)
$d.Content.Find.Execute("75×46=3450", $true, $false, $false, $false, $false, $true, 1, $false, "15×49=735", 2)
$d.Content.Find.Execute("93×55=5115", $true, $false, $false, $false, $false, $true, 1, $false, "93×44=4092", 2)

Write-Output "Done applying replacements"
